# EECE.3220 Lecture 37 (Exam 3 Preview) update
#  - Review: Heaps slide (slide 13): correct the array-indexing bullets
#  - All content slides (2-19): refresh the auto date placeholder text
#    from 12/10/2019 -> 12/11/2019

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 13 ("Review: Heaps") - fix the heap array-index bullets
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$bodyShape = $null
for ($i = 1; $i -le $slide13.Shapes.Count; $i++) {
    $sh = $slide13.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -like "*Children are at indexes*") {
        $bodyShape = $sh
        break
    }
}

if ($bodyShape -ne $null) {
    $tr = $bodyShape.TextFrame.TextRange

    # -- "If first element @ index 1, ..." -> "... index 0, ..."
    $t = $tr.Text
    $old = "If first element @ index 1, then, given index n"
    $new = "If first element @ index 0, then, given index n"
    $idx = $t.IndexOf($old)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $old.Length)
        $sub.Text = $new
    }

    # -- "Children are at indexes 2*n and 2*n + 1" -> "... 2*n + 1 and 2*n + 2"
    $t = $tr.Text
    $old = "Children are at indexes 2*n and 2*n + 1"
    $new = "Children are at indexes 2*n + 1 and 2*n + 2"
    $idx = $t.IndexOf($old)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $old.Length)
        $sub.Text = $new
    }

    # -- "Parent is at index n/2" -> "Parent is at " + "index (n-1)/2"
    #    (only the second half is rewritten, so the edit naturally lands in
    #    its own run, matching how the suffix was retyped by hand)
    $t = $tr.Text
    $old = "index n/2"
    $new = "index (n-1)/2"
    $idx = $t.IndexOf("Parent is at " + $old)
    if ($idx -ge 0) {
        $suffixStart = $idx + "Parent is at ".Length
        $sub = $tr.Characters($suffixStart + 1, $old.Length)
        $sub.Text = $new
    }
}

# ---------------------------------------------------------------------
# 2) Refresh the "datetime1" auto-date placeholder on every slide that
#    still shows the old capture date.
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tf = $sh.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                if ($tr.Text -eq "12/10/2019") {
                    $tr.Text = "12/11/2019"
                }
            }
        }
    }
}
